$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Exact "50×12=" "89×13="
Replace-Exact "27×42=" "41×49="
Replace-Exact "27×32=" "82×57="
Replace-Exact "28×21=" "54×31="
Replace-Exact "63×38=" "30×61="
Replace-Exact "60×45=" "53×18="
Replace-Exact "49×79=" "42×82="
Replace-Exact "34×40=" "56×27="
Replace-Exact "30×49=" "79×71="
Replace-Exact "74×47=" "74×77="
Replace-Exact "63×48=" "60×21="
Replace-Exact "14×35=" "16×26="
Replace-Exact "87×31=" "31×32="
Replace-Exact "29×41=" "44×59="
Replace-Exact "25×54=" "57×58="
Replace-Exact "67×32=" "96×13="
Replace-Exact "74×63=" "61×32="
Replace-Exact "71×38=" "97×60="
Replace-Exact "95×46=" "32×94="
Replace-Exact "17×97=" "22×38="
Replace-Exact "22×76=" "53×59="
Replace-Exact "35×50=" "43×59="
Replace-Exact "48×78=" "46×32="
Replace-Exact "40×41=" "61×14="
Replace-Exact "91×80=" "85×32="
